# Fruta / hortaliza, semanal
#
# A new weekly price-report record needs to be inserted into the
# "Feria Lagunitas de Puerto Montt - Membrillo" data block. The block runs
# from row 90 to row 164 (all rows share the same Mercado/Región/Producto/
# Variedad/Origen metadata in columns A,B,C,E,F,G,H,I,J,K,R); only the
# per-record columns D (Fecha), L (Calidad), M (Volumen), N (Precio
# mínimo), O (Precio máximo), P (Precio promedio ponderado), Q (Unidad de
# comercialización), S (Precio $/Kg) and T (Kg / unidad) vary.
#
# The new record is inserted at row 90, pushing every existing record in
# the block down by one row (old row 90 -> new row 91, ..., old row 164 ->
# new row 165), exactly like Excel's native "Insert Sheet Rows" command.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 90; this shifts rows 90:164 down to 91:165,
# carrying their values and formatting along (mirrors Excel's row insert).
$ws.Rows(90).Insert()

# Populate the newly-inserted row 90 with the new weekly record.
$ws.Range("A90").Value = 4
$ws.Range("B90").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C90").Value = "Los Lagos"
$ws.Range("D90").Value = 45090
$ws.Range("E90").Value = 10
$ws.Range("F90").Value = "Fruta"
$ws.Range("G90").Value = 100104
$ws.Range("H90").Value = "Frutos de pepita"
$ws.Range("I90").Value = 100104003
$ws.Range("J90").Value = "Membrillo"
$ws.Range("K90").Value = "Champion"
$ws.Range("L90").Value = "Primera"
$ws.Range("M90").Value = 300
$ws.Range("N90").Value = 12000
$ws.Range("O90").Value = 13000
$ws.Range("P90").Value = 12500
$ws.Range("Q90").Value = '$/caja 18 kilos empedrada'
$ws.Range("R90").Value = "Región de O'Higgins"
$ws.Range("S90").Value = 694
$ws.Range("T90").Value = 18
